# Team_02_M2_D1_Dimension_Mapping_Dealer.xlsx
# Re-aligns the "Dealer" sheet's Source-1 (D:F) / Source-2 (G:I) mapping
# columns for rows 4, 11 and 16 so the Corporate/Branch mappings line up
# with their correct source table, and updates the current selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dealer")

# --- Row 4 (DLR_Name row): swap the Source-1 / Source-2 table+column pair ---
$ws.Range("D4").Value = "SEIS732_Team_02_Corporate"
$ws.Range("E4").Value = "Dealership"
$ws.Range("G4").Value = "SEIS732_Team_02_Sales_Org"
$ws.Range("H4").Value = "Dealer"

# --- Row 11 (DLR_Zip row): swap the Source-1 / Source-2 table+column triple ---
$ws.Range("D11").Value = "SEIS732_Team_02_Corporate"
$ws.Range("E11").Value = "Branch"
$ws.Range("F11").Value = "BR_Zip"
$ws.Range("G11").Value = "SEIS732_Team_02_Sales_Org"
$ws.Range("H11").Value = "Dealer"
$ws.Range("I11").Value = "DLR_Zip"

# --- Row 16 (DLR_Phone row): swap the Source-1 / Source-2 table+column triple ---
$ws.Range("D16").Value = "SEIS732_Team_02_Corporate"
$ws.Range("E16").Value = "Branch"
$ws.Range("F16").Value = "BR_Phone"
$ws.Range("G16").Value = "SEIS732_Team_02_Sales_Org"
$ws.Range("H16").Value = "Dealer"
$ws.Range("I16").Value = "DLR_Phone"

# --- Update the sheet's saved selection to the header cells A1:C1 ---
$ws.Range("A1:C1").Select()
